# 03.11.2024 - PART 1
# Add IFNA / IFERROR rows to the FUNCTIONS reference sheet, and update the
# saved cursor/selection state on the FUNCTIONS and OTHER sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# FUNCTIONS sheet: append IFNA (row 4) and IFERROR (row 5), cloning the
# formatting of the existing IF row (row 3).
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("FUNCTIONS")

$ws.Range("A3:E3").Copy()
$ws.Range("A4:E5").PasteSpecial(-4122)   # xlPasteFormats

# Column A / C: EXCEL_FUNCTION / PYTHON_FUNCTION name
$ws.Cells.Item(4, 1).Value = "IFNA"
$ws.Cells.Item(5, 1).Value = "IFERROR"

# Column B: JSON_FIELDS
$ws.Cells.Item(4, 2).Value = "['CONDITION', 'CORRECTION']"
$ws.Cells.Item(5, 2).Value = "['CONDITION', 'CORRECTION']"

# Column C: PYTHON_FUNCTION
$ws.Cells.Item(4, 3).Value = "IFNA"
$ws.Cells.Item(5, 3).Value = "IFERROR"

# Column D: PYTHON_FILE
$ws.Cells.Item(4, 4).Value = "EXCEL_FUNCTIONS"
$ws.Cells.Item(5, 4).Value = "EXCEL_FUNCTIONS"

# Column E: CLUSTERS
$ws.Cells.Item(4, 5).Value = "['CONDITION']"
$ws.Cells.Item(5, 5).Value = "['CONDITION']"

# ---------------------------------------------------------------------
# OTHER sheet: move the saved selection off the "select-all" state onto
# cell B24.
# ---------------------------------------------------------------------
$wsOther = $wb.Worksheets.Item("OTHER")
$wsOther.Range("B24").Select()

# ---------------------------------------------------------------------
# Re-activate FUNCTIONS and move its saved selection to B9 (also drops
# the stale topLeftCell scroll position from the old selection).
# ---------------------------------------------------------------------
$ws.Activate()
$ws.Range("B9").Select()
